# Generate Report for Handoff
#
# Semantic changes (derived from the canonical OOXML diff):
#  1. Priority column ("low" -> "ht") for the handed-off files on both the
#     zh-cn and de-de sheets (rows 4-7, column E).
#  2. Latest Handoff Datetime for zh-cn (column H, rows 4-7):
#       2016-08-27 10:31:41 -> 2016-08-27 10:32:13
#  3. Latest HO Xliff Generate Date, which is shared between the Overview
#     sheet (column G, rows 4-7) and the de-de sheet (column H, rows 4-7):
#       2016-08-27 10:31:47 -> 2016-08-27 10:32:17

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$rows = @(4, 5, 6, 7)

foreach ($r in $rows) {
    # Priority: low -> ht (both language sheets)
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"

    # zh-cn: Latest Handoff Datetime refreshed
    $wsZhCn.Range("H$r").Value = "2016-08-27 10:32:13"

    # Latest HO Xliff Generate Date refreshed - shows up both on the
    # Overview sheet and on the de-de sheet
    $wsOverview.Range("G$r").Value = "2016-08-27 10:32:17"
    $wsDeDe.Range("H$r").Value = "2016-08-27 10:32:17"
}
